$d = $word.ActiveDocument

# The document ends with (in order, right before the sectPr):
#   ... "Aplicaciones"
#   <empty paragraph, pStyle p2>   <- stays as-is
#   <empty paragraph, pStyle p2>   <- second-to-last paragraph; remove it
#   <empty paragraph, no style>    <- last paragraph; becomes the new
#                                     "Referencias" paragraph
#
# Deleting the Range of the second-to-last (empty, "p2"-styled) paragraph
# merges it away, leaving the final paragraph mark intact.
$toRemove = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$toRemove.Range.Delete()

# The (now) last paragraph in the document becomes the "Referencias" one.
$last = $d.Paragraphs.Last

$text = "Referencias: Ash (Basic Probability Theory, Probability and Measure Theory, Real Analysis and Probability), Billingsley, Ross(A first course in probability, Introduction ton probability models-sexta,novena y décima edición), Mood, Casella (Statistical Inference), Papoulis (Probability, Random Variables and stochastic processes), Gubner (Probability and Random Processes for Electrical and Computer Engineers), Kurtz (Probability and its applications),  Rincón (Curso Elemental de Probabilidad, Curso Intermedio de Probabilidad), Grinstead (Introduction to Probability), Feller (Introduction to Probability Theory and its applications), Kay (Intuitive Probability and Random Processes using Matlab), Pitman (Probability)."

$last.Range.InsertAfter($text)

# Give the new run (and the paragraph mark) the small 6.5pt (half-point
# size 13) font used for the references line, matching both <w:sz> and
# <w:szCs>.
$fmtRange = $last.Range
$fmtRange.Font.Size = 6.5
$fmtRange.Font.SizeBi = 6.5

Write-Output "done"
